# Lecture 4/5 update:
#  1) bump the stale cached "last saved" datetimeFigureOut fields from
#     1/25/18 to 1/26/18, wherever one happens to live on this deck; and
#  2) fix the subscripted-variable labels in the k-mer/discretization
#     diagram (Z_m, V_m, X_n, Y_n -> z_i, y_i, s_i, u_i).

$p = $ppt.ActivePresentation

# --- 1) date field refresh (defensive: only touches runs that actually
#        contain the stale cached date string, so it's a no-op on slides
#        without a datetimeFigureOut field) ---------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -like "*1/25/18*") {
                $runCount = $tr.Runs().Count
                for ($ri = 1; $ri -le $runCount; $ri++) {
                    $run = $tr.Runs($ri)
                    if ($run.Text -like "*1/25/18*") {
                        $run.Text = $run.Text.Replace("1/25/18", "1/26/18")
                    }
                }
            }
        }
    }
}

# --- 2) k-mer/discretization diagram labels --------------------------------
$s = $p.Slides.Item(1)

# Map of shape index -> (base-letter text, subscript-letter text)
$fixes = @{
    5  = @("z", "i")   # TextBox 12: "Z" / "m" (baseline -25000)  -> "z" / "i"
    7  = @("y", "i")   # TextBox 16: "V" / "m" (baseline -25000)  -> "y" / "i"
    10 = @("s", "i")   # TextBox 20: "X" / "n" (baseline -25000)  -> "s" / "i"
    11 = @("u", "i")   # TextBox 21: "Y" / "n" (baseline -25000)  -> "u" / "i"
}

foreach ($idx in $fixes.Keys) {
    $shape = $s.Shapes.Item($idx)
    $tr = $shape.TextFrame.TextRange
    $count = $tr.Runs().Count
    # last two runs are the italic base letter and its subscripted index letter
    $tr.Runs($count - 1).Text = $fixes[$idx][0]
    $tr.Runs($count).Text = $fixes[$idx][1]
}
